$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.831.71'
$ws.Range("E2").Value = '  +11.74%  '
$ws.Range("D3").Value = '1.736.32'
$ws.Range("E3").Value = '  +7.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").Value = '  +1.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.27'
$ws.Range("E5").Value = '  +4.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9936'
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3821'
$ws.Range("E7").Value = '  +4.57%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3646'
$ws.Range("E8").Value = '  +6.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.45'
$ws.Range("E9").Value = '  +18.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.226'
$ws.Range("E10").Value = '  +7.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07707'
$ws.Range("E11").Value = '  +9.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9934'
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.83'
$ws.Range("E13").Value = '  +9.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.477'
$ws.Range("E14").Value = '  +10.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.109'
$ws.Range("E15").Value = '  +7.75%  '
$ws.Range("D16").Value = '1.759.45'
$ws.Range("E16").Value = '  +9.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001155'
$ws.Range("E17").Value = '  +6.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9925'
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06814'
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '87.08'
$ws.Range("E20").Value = '  +11.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.52'
$ws.Range("E21").Value = '  +8.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.482'
$ws.Range("E22").Value = '  +8.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.77'
$ws.Range("E23").Value = '  +9.40%  '
$ws.Range("D24").Value = '25.730.19'
$ws.Range("E24").Value = '  +11.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.438'
$ws.Range("E25").Value = '  +3.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.929'
$ws.Range("E26").Value = '  +12.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.58'
$ws.Range("E27").Value = '  +6.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '154.94'
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '134.36'
$ws.Range("E29").Value = '  +7.88%  '
$ws.Range("D30").Value = '1.927.99'
$ws.Range("E30").Value = '  +8.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.197'
$ws.Range("E31").Value = '  +22.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.040'
$ws.Range("E32").Value = '  +16.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.183'
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.23'
$ws.Range("E34").Value = '  +20.20%  '
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.791'
$ws.Range("E35").Value = '  +7.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08680'
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.648'
$ws.Range("E37").Value = '  +9.37%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06751'
$ws.Range("E38").Value = '  +10.59%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02476'
$ws.Range("E39").Value = '  +11.88%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.274'
$ws.Range("E40").Value = '  +8.38%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.302'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2220'
$ws.Range("E42").Value = '  +10.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6518'
$ws.Range("E43").Value = '  +10.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9927'
$ws.Range("E44").Value = '  +1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.95'
$ws.Range("E45").Value = '  +7.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6326'
$ws.Range("E46").Value = '  +10.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.905'
$ws.Range("E47").Value = '  +3.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.176'
$ws.Range("E48").Value = '  +10.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.41'
$ws.Range("E49").Value = '  +4.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07443'
$ws.Range("E50").Value = '  +7.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.53'
$ws.Range("E51").Value = '  +8.32%  '
